$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.930.63'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '1.904.85'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5030'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4039'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08254'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.97'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.094'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.12%  '
$ws.Range('E12').Value = '  +2.39%  '
$ws.Range('D13').Value = '1.904.29'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.353'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.174'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001092'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06490'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.53%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.931'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').Value = '29.971.85'
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.190'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '22.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.26%  '
$ws.Range('D27').Value = '2.123.86'
$ws.Range('E27').Value = '  -0.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.263'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.51'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.120'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.47%  '
$ws.Range('E32').Value = '  -2.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.893'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.795'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02428'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.73%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.362'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.28%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06326'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2138'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.10%  '
$ws.Range('E39').Value = '  -3.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6421'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.600'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.30'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.59%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.204'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.192'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.23'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5992'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.01%  '
$ws.Range('E47').Value = '  -2.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.18'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.203'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '78.31'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.127'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.08%  '
